# Sort each contiguous "year" block of the bank-rate table ascending by
# column A (the count column), keeping each row's E value attached to its
# own A value (B/C/D are constant within a block, so they are unaffected
# either way).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row ranges (inclusive) for each contiguous "year" block in the data.
$groups = @(
    @(2, 5),
    @(6, 10),
    @(11, 15),
    @(16, 20),
    @(21, 22)
)

foreach ($grp in $groups) {
    $startRow = $grp[0]
    $endRow = $grp[1]
    $n = $endRow - $startRow + 1

    # Capture the original A and E values for the block (Value2 for
    # reliable numeric reads).
    $aVals = @()
    $eVals = @()
    for ($r = $startRow; $r -le $endRow; $r++) {
        $aVals += $ws.Cells.Item($r, 1).Value2
        $eVals += $ws.Cells.Item($r, 5).Value2
    }

    # Simple in-place insertion sort on the parallel A/E arrays, ascending
    # by A (Sort-Object with a key scriptblock proved unreliable on the
    # nested-array pairs in this runtime, so do it manually).
    for ($i = 1; $i -lt $n; $i++) {
        $keyA = $aVals[$i]
        $keyE = $eVals[$i]
        $j = $i - 1
        while ($j -ge 0 -and $aVals[$j] -gt $keyA) {
            $aVals[$j + 1] = $aVals[$j]
            $eVals[$j + 1] = $eVals[$j]
            $j--
        }
        $aVals[$j + 1] = $keyA
        $eVals[$j + 1] = $keyE
    }

    # Write the sorted values back into the block.
    for ($i = 0; $i -lt $n; $i++) {
        $targetRow = $startRow + $i
        $ws.Cells.Item($targetRow, 1).Value = $aVals[$i]
        $ws.Cells.Item($targetRow, 5).Value = $eVals[$i]
    }
}
